# Add a "2022-Q3" sheet (new first quarter in the series) right after
# "总计" and before "2022-Q2", fill it with its fund-holding table, and
# insert a matching summary row at the top of "总计"'s data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new summary row (old row 2 onward shift down by one) in
#    the "总计" sheet, and populate it with the 2022-Q3 totals.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows("2:2").Insert()

# The freshly inserted row inherited the header's style; re-stamp it with
# the same look used by every other data row (copy formats from the row
# that used to be row 2 and is now row 3).
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 10
$summary.Range("D2").Value = 1.2

# ---------------------------------------------------------------------
# 2) Create the new "2022-Q3" worksheet, positioned right before the
#    existing "2022-Q2" tab, and re-select the sheet that was active
#    before (Add() makes the new sheet active otherwise).
# ---------------------------------------------------------------------
$previouslyActiveName = $wb.ActiveSheet.Name
$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($beforeSheet)
$q3.Name = "2022-Q3"

# Header row.
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Bold + thin border + centered, matching the header/index styling used
# throughout the workbook; then fan that format out to the rest of the
# header row and to the index column below.
$q3.Range("B1").Font.Bold = $true
$q3.Range("B1").HorizontalAlignment = -4108
$q3.Range("B1").VerticalAlignment = -4160
$q3.Range("B1").Borders.LineStyle = 1
$q3.Range("B1").Copy()
$q3.Range("C1:H1").PasteSpecial(-4122)
$q3.Range("A2:A11").PasteSpecial(-4122)

# Data rows. Columns D/E/F/G hold numeric-looking figures that are
# actually stored as text (matching the source sheet), and column B's
# fund codes must stay text too (leading zeros). Using Formula with a
# leading apostrophe forces the text type without Excel re-parsing the
# numeric-looking string into a number.
function Set-Text($cell, $text) {
    $cell.Formula = "'" + $text
}

$rows = @(
    @("870009", "广发资管平衡精选一年持有混合A", "7.47", "92.14", "9.83", "0.7343", 1),
    @("014062", "景顺长城专精特新量化优选股票A", "8.02", "91.10", "1.79", "0.1436", 9),
    @("872019", "广发资管平衡精选一年持有混合C", "1.09", "92.14", "9.83", "0.1071", 1),
    @("014063", "景顺长城专精特新量化优选股票C", "5.41", "91.10", "1.79", "0.0968", 9),
    @("001917", "招商量化精选股票A", "4.16", "92.70", "1.65", "0.0686", 2),
    @("007950", "招商量化精选股票C", "2.39", "92.70", "1.65", "0.0394", 2),
    @("001375", "金元顺安优质精选灵活配置混合C", "0.62", "65.13", "0.74", "0.0046", 7),
    @("002952", "建信多因子量化股票", "0.09", "91.26", "3.43", "0.0031", 5),
    @("015245", "南华丰汇混合", "0.09", "86.53", "1.03", "0.0009", 9),
    @("620007", "金元顺安优质精选灵活配置混合A", "0.06", "65.13", "0.74", "0.0004", 7)
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $r - 2
    Set-Text $q3.Cells.Item($r, 2) $row[0]
    $q3.Cells.Item($r, 3).Value = $row[1]
    Set-Text $q3.Cells.Item($r, 4) $row[2]
    Set-Text $q3.Cells.Item($r, 5) $row[3]
    Set-Text $q3.Cells.Item($r, 6) $row[4]
    Set-Text $q3.Cells.Item($r, 7) $row[5]
    $q3.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# Restore whichever sheet was selected before we added the new tab, so
# the workbook's "active sheet" bookkeeping is left exactly as it was.
# (Looked up by name again - sheet references captured before Add()
# shift to track sheet *position*, not the sheet they originally
# pointed to.)
$wb.Worksheets.Item($previouslyActiveName).Activate()
